$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Highlight B8:K10 with the existing yellow fill style (same as used
#        elsewhere in the workbook, e.g. C1!B19) so that we reuse style index 7
#        rather than create a brand-new one. ---
$wsC1 = $wb.Worksheets.Item("C1")
$wsC1.Range("B19").Copy()
$ws.Range("B8:K10").PasteSpecial(-4122)

# --- 2. Row 16: AVERAGE of rows 8:10 per column (B..K), written as a shared
#        formula for C16:K18 (matching the original author's fill-down extent,
#        which still covered rows 17-18 before those extra cells were cleared). ---
$ws.Range("B16").Formula = "=AVERAGE(B8:B10)"
$ws.Range("C16:K18").Formula = "=AVERAGE(C8:C10)"
$ws.Range("C17:K18").ClearContents()

# --- 3. Rows 18-27: column B repeats the row-16 averages, column C expresses
#        them as a percentage of the first entry (style 1 = light blue, copied
#        from an existing cell using that style, e.g. C2!C26). ---
$ws.Range("B18").Value2 = 0.042333333333333334
$ws.Range("B19").Value2 = 0.044333333333333336
$ws.Range("B20").Value2 = 0.045333333333333337
$ws.Range("B21").Value2 = 0.044666666666666667
$ws.Range("B22").Value2 = 0.042333333333333334
$ws.Range("B23").Value2 = 0.04766666666666667
$ws.Range("B24").Value2 = 0.047333333333333338
$ws.Range("B25").Value2 = 0.046666666666666669
$ws.Range("B26").Value2 = 0.048333333333333339
$ws.Range("B27").Value2 = 0.048999999999999995

$ws.Range("C18").Value2 = 100
$ws.Range("C19").Formula = "=B19/0.0433*100"
$ws.Range("C20:C27").Formula = "=B20/0.0433*100"

$wsC2 = $wb.Worksheets.Item("C2")
$wsC2.Range("C26").Copy()
$ws.Range("C18:C27").PasteSpecial(-4122)

# --- 4. Update the selection to match the saved view state. ---
$ws.Range("H23").Select()
